# Insert two new rows before existing row 203, shifting the rest of the
# data block (old rows 203-249) down to rows 205-251, then populate the
# two new rows with their own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 203 (existing rows 203.. shift down by 2)
$ws.Rows("203:204").Insert()

# New row 203: Choclero / Primera
$ws.Range("A203").Value = 4
$ws.Range("B203").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C203").Value = "Los Lagos"
$ws.Range("D203").Value = "2022-03-17"
$ws.Range("D203").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 100112024
$ws.Range("G203").Value = "Choclo"
$ws.Range("H203").Value = "Choclero"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 3000
$ws.Range("K203").Value = 250
$ws.Range("L203").Value = 250
$ws.Range("M203").Value = 250
$ws.Range("N203").Value = "`$/unidad"
$ws.Range("O203").Value = "Región del Maule"
$ws.Range("P203").Value = 250
$ws.Range("Q203").Value = 1
$ws.Range("R203").Value = "Hortaliza"

# New row 204: Dulce o Americano / Primera
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = "2022-03-17"
$ws.Range("D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = 100112024
$ws.Range("G204").Value = "Choclo"
$ws.Range("H204").Value = "Dulce o Americano"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 6000
$ws.Range("K204").Value = 180
$ws.Range("L204").Value = 200
$ws.Range("M204").Value = 190
$ws.Range("N204").Value = "`$/unidad"
$ws.Range("O204").Value = "Región de Los Lagos"
$ws.Range("P204").Value = 190
$ws.Range("Q204").Value = 1
$ws.Range("R204").Value = "Hortaliza"
